# Parsiss TimeSheet update: add "اسفند 98" (Esfand 98) month section with
# a "Multithreading" activity ("* Multithread Tracking" task), and rename
# two existing task labels ("* Virtual Camera" -> "* 3D Virtual View",
# "* Apply Registration" -> "* Live Tracking (Pre-registered)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename two existing task labels in the "آبان تا بهمن 98" block ----
$ws.Range("B67").Value = "* 3D Virtual View"
$ws.Range("E66").Value = "* Live Tracking (Pre-registered)"

# --- 2. Fill in the content of the new block first (values/formulas) ------
# Doing this BEFORE copying formatting avoids a quote-prefix/style quirk
# that appears when a string value is assigned after the format is applied.

# Header row (month name / column headers)
$ws.Range("A77").Value = "اسفند 98"
$ws.Range("B77").Value = "Activity"
$ws.Range("C77").Value = "Hours"
$ws.Range("E77").Value = "Tasks Done"

# First (and only populated) task row
$ws.Range("B78").Value = "Multithreading"
$ws.Range("C78").Value = 1
$ws.Range("E78").Value = "* Multithread Tracking"

# Totals
$ws.Range("B87").Value = "• Total Hours"
$ws.Range("C87").Formula = "=SUM(C78:C86)"

$ws.Range("C88").Value = "@Parsiss"
$ws.Range("D88").Value = 0

$ws.Range("C89").Value = "@Home"
$ws.Range("D89").Value = 1

# --- 3. Apply formatting, mirroring rows 60-74 onto rows 75-89 ------------
# Copy formatting only (no values) cell-by-cell from the existing block so
# the new block picks up identical styles without disturbing used-range.

$fmtPairs = @(
    @("A60","A75"), @("B60","B75"), @("C60","C75"), @("D60","D75"),
    @("E60","E75"), @("F60","F75"), @("G60","G75"),

    @("A61","A76"), @("B61","B76"), @("C61","C76"), @("D61","D76"), @("E61","E76"),

    @("A62","A77"), @("B62","B77"), @("C62","C77"), @("E62","E77"),

    @("B63","B78"), @("C63","C78"), @("E63","E78"),

    @("B64","B79"), @("C64","C79"), @("E64","E79"),

    @("B65","B80"), @("C65","C80"), @("E65","E80"),

    @("B66","B81"), @("C66","C81"), @("E66","E81"),

    @("B67","B82"), @("C67","C82"),
    @("B68","B83"), @("C68","C83"),
    @("B69","B84"), @("C69","C84"),
    @("B70","B85"), @("C70","C85"),
    @("B71","B86"), @("C71","C86"),

    @("B72","B87"), @("C72","C87"),

    @("C73","C88"), @("D73","D88"),

    @("C74","C89"), @("D74","D89")
)

foreach ($pair in $fmtPairs) {
    $srcAddr = $pair[0]
    $dstAddr = $pair[1]
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- 4. Keep the view selection in sync with the new bottom of the sheet --
$ws.Range("E86").Select()
